$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old range (A1:F5) so leftover cells (C-F, rows 3-5) are removed
$ws.Range("A1:F5").Clear()

# Write the new values
$ws.Range("A1").Value = 44
$ws.Range("B1").Value = 22
$ws.Range("C1").Value = 1
$ws.Range("A2").Value = 48
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 1

# Update the selection to match the target state
$ws.Range("B4").Select()
